$d = $word.ActiveDocument

# Locate the end of the "Disciplinary action form" title run so we can
# remove the trailing lone-space run that follows it and drop a
# zero-width "_GoBack" bookmark in its place (this also relocates the
# existing "_GoBack" bookmark away from its old spot further down the
# document, since Word only allows one bookmark per name).
$titleRange = $d.Content
$found = $titleRange.Find.Execute("Disciplinary action form", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$titleEnd = $titleRange.End

# The run immediately after the title is a single space character -
# delete just that character, leaving the title run's own formatting
# (and its rsid attributes) untouched.
$spaceRange = $d.Range($titleEnd, $titleEnd + 1)
$spaceRange.Delete()

# Re-add (relocate) the "_GoBack" bookmark as a zero-length bookmark
# right after the title text, before the following line-break run.
$bmRange = $d.Range($titleEnd, $titleEnd)
$d.Bookmarks.Add("_GoBack", $bmRange)
